$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.928.59'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '1.637.81'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.60%  '
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5075'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06367'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.78'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07761'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.297'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').Value = '1.639.55'
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5459'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('D15').Value = '0.0₅7742'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.16'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').Value = '25.940.27'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.461'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '196.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.945'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.146'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.890'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1249'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.840'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.237'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04883'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.247'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.205'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.371'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9137'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.570'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5525'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('D38').Value = '1.122.62'
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01566'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.001'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.604'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.54%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8038'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '98.52'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.58%  '
$ws.Range('D44').Value = '0.0₈121'
$ws.Range('E44').Value = '  -9.38%  '
$ws.Range('D45').Value = '1.774.29'
$ws.Range('E45').Value = '  -0.39%  '
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.14'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('E48').Value = '  -0.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05180'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.492'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.002'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.55%  '
